$d = $word.ActiveDocument

# Locate the "b. Faltaron atributos." paragraph inside section "2. MC:" --
# its predecessor ("a. No se realizó...") and successor ("c. Algunas
# navegabilidades...") form the three still-unstruck items of that list
# that need to be marked as resolved (strikethrough).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Faltaron atributos.*") {
        $p.Previous(1).Range.Font.StrikeThrough = 1
        $p.Range.Font.StrikeThrough = 1
        $p.Next(1).Range.Font.StrikeThrough = 1
        Write-Output "Struck MC a/b/c"
    }
    if ($p.Range.Text -like "*La clase Pedido necesita de otras para mapear objetos.*") {
        $p.Range.Font.StrikeThrough = 1
        Write-Output "Struck Capas d"
    }
}
